# Applies the 27-11-2023 14:45 scraper update to the Indonesia Liga-1
# 2023-2024 sheet:
#   1) Thirteen pairs of adjacent rows had their match data (columns F:V)
#      swapped between them (the row's Indice/pais/torneio/temporada/
#      data_partida in A:E stayed put; only the match payload moved).
#   2) Two brand-new match rows (178 and 179 in 1-based sheet terms, i.e.
#      spreadsheet rows 177/178 of data) were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row-pair swaps (columns F:V only) -----------------------------
$swapPairs = @(
  @(17, 18),
  @(19, 20),
  @(30, 31),
  @(42, 43),
  @(44, 45),
  @(82, 83),
  @(84, 85),
  @(93, 94),
  @(103, 104),
  @(107, 108),
  @(134, 135),
  @(162, 163),
  @(165, 166)
)

foreach ($pair in $swapPairs) {
  $rowA = $pair[0]
  $rowB = $pair[1]
  $rangeA = $ws.Range("F$rowA" + ":V$rowA")
  $rangeB = $ws.Range("F$rowB" + ":V$rowB")
  $valuesA = $rangeA.Value2
  $valuesB = $rangeB.Value2
  $rangeA.Value2 = $valuesB
  $rangeB.Value2 = $valuesA
}

# --- 2) Append two new match rows --------------------------------------
$newRows = @(
  @{
    Row = 177
    Index = 176
    Date = 45257.54166666666
    Home = "Borneo"
    HomeGoals = 1
    Away = "Persis Solo"
    AwayGoals = 0
    HomeOpenOdds = 1.56
    HomeOpenTime = "26/11/2023 01:12"
    HomeCloseOdds = 1.47
    HomeCloseTime = "27/11/2023 12:55"
    DrawOpenOdds = 4.02
    DrawOpenTime = "26/11/2023 01:12"
    DrawCloseOdds = 4.28
    DrawCloseTime = "27/11/2023 12:58"
    AwayOpenOdds = 4.71
    AwayOpenTime = "26/11/2023 01:12"
    AwayCloseOdds = 6.96
    AwayCloseTime = "27/11/2023 12:55"
    Url = "https://www.betexplorer.com/football/indonesia/liga-1/borneo-persis-solo/lvp3JeAj/"
  },
  @{
    Row = 178
    Index = 177
    Date = 45257.54166666666
    Home = "FC Bhayangkara"
    HomeGoals = 2
    Away = "Persija Jakarta"
    AwayGoals = 2
    HomeOpenOdds = 3.03
    HomeOpenTime = "26/11/2023 01:12"
    HomeCloseOdds = 3.94
    HomeCloseTime = "27/11/2023 12:59"
    DrawOpenOdds = 3.3
    DrawOpenTime = "26/11/2023 01:12"
    DrawCloseOdds = 3.52
    DrawCloseTime = "27/11/2023 12:59"
    AwayOpenOdds = 2.14
    AwayOpenTime = "26/11/2023 01:12"
    AwayCloseOdds = 1.93
    AwayCloseTime = "27/11/2023 12:57"
    Url = "https://www.betexplorer.com/football/indonesia/liga-1/fc-bhayangkara-persija-jakarta/Gpxksa2c/"
  }
)

foreach ($row in $newRows) {
  $r = $row.Row

  $ws.Range("A$r").Value2 = $row.Index
  $ws.Range("B$r").Value2 = "indonesia"
  $ws.Range("C$r").Value2 = "liga-1"
  $ws.Range("D$r").Value2 = "2023-2024"
  $ws.Range("E$r").Value2 = $row.Date
  $ws.Range("F$r").Value2 = $row.Home
  $ws.Range("G$r").Value2 = $row.HomeGoals
  $ws.Range("H$r").Value2 = $row.Away
  $ws.Range("I$r").Value2 = $row.AwayGoals
  $ws.Range("J$r").Value2 = $row.HomeOpenOdds
  $ws.Range("K$r").Value2 = $row.HomeOpenTime
  $ws.Range("L$r").Value2 = $row.HomeCloseOdds
  $ws.Range("M$r").Value2 = $row.HomeCloseTime
  $ws.Range("N$r").Value2 = $row.DrawOpenOdds
  $ws.Range("O$r").Value2 = $row.DrawOpenTime
  $ws.Range("P$r").Value2 = $row.DrawCloseOdds
  $ws.Range("Q$r").Value2 = $row.DrawCloseTime
  $ws.Range("R$r").Value2 = $row.AwayOpenOdds
  $ws.Range("S$r").Value2 = $row.AwayOpenTime
  $ws.Range("T$r").Value2 = $row.AwayCloseOdds
  $ws.Range("U$r").Value2 = $row.AwayCloseTime
  $ws.Range("V$r").Value2 = $row.Url

  # Carry over the same cell formatting the rest of the table uses:
  # bold/centered/bordered for column A, the datetime number format for
  # column E (copied from the row right above, which still has the
  # original styling untouched by the swaps above).
  $srcRow = $r - 1
  $ws.Range("A$srcRow").Copy() | Out-Null
  $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
  $ws.Range("E$srcRow").Copy() | Out-Null
  $ws.Range("E$r").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false
